$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, shifting rows 96:101 down to 97:102.
# Excel copies the formatting of the row above (row 95) by default when
# inserting via Rows.Insert(), but since row 96's original content is what
# we want replicated (minus a handful of changed fields), copy row 96's
# values into the freshly inserted row explicitly.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the same "template" values that
# used to live in (old) row 96 / (new) row 97, since the new record is
# another Arandano (blue) / Macroferia Regional de Talca entry.
$ws.Cells.Item(96, 1).Value = 5
$ws.Cells.Item(96, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(96, 3).Value = "Maule"
$ws.Cells.Item(96, 4).Value = 44931
$ws.Cells.Item(96, 5).Value = 7
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100101
$ws.Cells.Item(96, 8).Value = "Berries"
$ws.Cells.Item(96, 9).Value = 100101001
$ws.Cells.Item(96, 10).Value = "Arándano (blue)"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 150
$ws.Cells.Item(96, 14).Value = 3000
$ws.Cells.Item(96, 15).Value = 3000
$ws.Cells.Item(96, 16).Value = 3000
$ws.Cells.Item(96, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(96, 19).Value = 1500
$ws.Cells.Item(96, 20).Value = 2
